# Natmi following Dr Hou advice
# Update recomputed ligand/receptor expression statistics for rows 2-4
# on the active sheet, reflecting a change in Ligand/Receptor-expressing
# cell counts (1 -> 3) and the resulting cascade of derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{
        E = 3
        G = 5.084811666666667
        H = 15.254435
        I = 0.5108142255059609
        J = 0.5108142255059609
        K = 3
        M = 0.8412133333333333
        N = 2.52364
        Q = 4.277411371488888
        R = 38.4967023434
        S = 0.5108142255059609
        T = 0.5108142255059609
    }
    3 = @{
        E = 3
        G = 2.689016666666667
        H = 8.06705
        I = 0.2701354653822224
        J = 0.2701354653822224
        K = 3
        M = 0.8412133333333333
        N = 2.52364
        Q = 2.262036673555555
        R = 20.358330062
        S = 0.2701354653822224
        T = 0.2701354653822224
    }
    4 = @{
        E = 3
        G = 2.180498333333333
        H = 6.541495
        I = 0.2190503091118167
        J = 0.2190503091118167
        K = 3
        M = 0.8412133333333333
        N = 2.52364
        Q = 1.834264271311111
        R = 16.5083784418
        S = 0.2190503091118167
        T = 0.2190503091118167
    }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$rowNum"
        $ws.Range($cellRef).Value = $cols[$col]
    }
}
